$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "il ha gros <m>os</m>" -> "il ha <m>gros os</m>"
#   - run "il ha gros " loses the trailing "gros "
#   - run "os" (inside the blue <m>...</m> markup) gains a "gros " prefix
# ---------------------------------------------------------------------------
$r1 = $d.Content
$null = $r1.Find.Execute("il ha gros ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$grosEnd = $r1.End
$grosRange = $d.Range($grosEnd - 5, $grosEnd)
$grosRange.Text = ""

$r1b = $d.Content
$null = $r1b.Find.Execute("os</m>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$osRange = $d.Range($r1b.Start, $r1b.Start + 2)
$osRange.Text = "gros os"

# ---------------------------------------------------------------------------
# Change 2: the <m> markup moves from before "moules" to before "lopins"
#   "... & les lopins aussy<lb/>des <m>moules qui ont servy</m> ..."
#     -> "... & les <m>lopins aussy<lb/>des moules qui ont servy</m> ..."
# ---------------------------------------------------------------------------

# 2a. remove the existing "<m>" that sits right before "moules qui ont servy"
$r2 = $d.Content
$null = $r2.Find.Execute("des <m>moules qui ont servy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mOldStart = $r2.Start + 4
$mOldRange = $d.Range($mOldStart, $mOldStart + 3)
$mOldRange.Text = ""

# 2b. insert a new "<m>" right before "lopins aussy", styled like the other
#     <m> markup runs (Courier New, blue, 9pt)
$r3 = $d.Content
$null = $r3.Find.Execute(" les lopins aussy", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$lopinsStart = $r3.Start + 5
$insPoint = $d.Range($lopinsStart, $lopinsStart)
$insPoint.InsertBefore("<m>")

$mNewRange = $d.Range($lopinsStart, $lopinsStart + 3)
$mNewRange.Font.Name = "Courier New"
$mNewRange.Font.Color = 16711680
$mNewRange.Font.Size = 9

# ---------------------------------------------------------------------------
# Change 3: "<m>ardille</m></tl> fresche. " -> "<m>ardille fresche</m></tl>. "
#   "fresche" moves from after the closing </m></tl> markup to before it,
#   staying adjacent to "ardille".
# ---------------------------------------------------------------------------
$r4 = $d.Content
$null = $r4.Find.Execute("ardille", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r4.Text = "ardille fresche"

$r5 = $d.Content
$null = $r5.Find.Execute(" fresche. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r5.Text = ". "
